$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 97.87634702449229
$ws.Range("H2").Value = 97.54939900844225
$ws.Range("I2").Value = 96.29110822582112

$ws.Range("G3").Value = 98.18353953819677
$ws.Range("H3").Value = 97.63649770943671
$ws.Range("I3").Value = 96.23228368043365

$ws.Range("G4").Value = 98.04254249801228
$ws.Range("H4").Value = 97.47168335959722
$ws.Range("I4").Value = 96.24212606346477

$ws.Range("G5").Value = 97.95454774210143
$ws.Range("H5").Value = 97.51271022009249
$ws.Range("I5").Value = 96.21058481004566

$ws.Range("G6").Value = 98.0669061578936
$ws.Range("H6").Value = 97.53833081934947
$ws.Range("I6").Value = 96.13680781001516
